$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 276, shifting existing rows 276:370 down to 277:371
$ws.Rows.Item(276).Insert()

# Populate the new row 276 with the new record's data.
$ws.Cells.Item(276, 1).Value = 3
$ws.Cells.Item(276, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(276, 3).Value = "Coquimbo"
$ws.Cells.Item(276, 4).Value = 44524
$ws.Cells.Item(276, 5).Value = 5
$ws.Cells.Item(276, 6).Value = 100112045
$ws.Cells.Item(276, 7).Value = "Zapallo"
$ws.Cells.Item(276, 8).Value = "Camote"
$ws.Cells.Item(276, 9).Value = "1a nueva(o)"
$ws.Cells.Item(276, 10).Value = 120
$ws.Cells.Item(276, 11).Value = 600
$ws.Cells.Item(276, 12).Value = 600
$ws.Cells.Item(276, 13).Value = 600
$ws.Cells.Item(276, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(276, 15).Value = "Provincia de Talca"
$ws.Cells.Item(276, 16).Value = 600
$ws.Cells.Item(276, 17).Value = 1
$ws.Cells.Item(276, 18).Value = "Hortaliza"
